$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "23.03.2023 10:08 (CET)"
$ws.Range("C4").Value = "https://gitlab.intra.infineon.com/semantic-web-projects/digital-reference/process_model_version/-/commit/69c451ed3cc43181ccf7143ba7e6ada345d38d6f"
$ws.Range("D4").Value = "7d94edb76e6b01e30e2dea18580509d0f3efcd866bfcad2a4b298e040d1da0f6"
